# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Nuevos periodos de mora (columna E), ahora en orden ascendente 2302..2308
$ws.Range("E16").Value = "2302"
$ws.Range("E17").Value = "2303"
$ws.Range("E18").Value = "2304"
$ws.Range("E19").Value = "2305"
$ws.Range("E20").Value = "2306"
$ws.Range("E21").Value = "2307"
$ws.Range("E22").Value = "2308"

# Valor Mora (columna F): el periodo parcial (46400 -> 40000) ahora cae en la
# ultima fila (2308) en vez de la primera (2302)
$ws.Range("F16").Value = 46400
$ws.Range("F17").Value = 46400
$ws.Range("F18").Value = 46400
$ws.Range("F19").Value = 46400
$ws.Range("F20").Value = 46400
$ws.Range("F21").Value = 46400
$ws.Range("F22").Value = 40000

# Salario Basico (columna G): actualizado de 1160000 a 1000000 para todas las filas
$ws.Range("G16:G22").Value = 1000000
